# Commit: "Used the serializable dictionary instead of list to store skill config"
#
# The skill config record was renamed from the list-backed "skill.TbSkill"
# to the dictionary-backed "TbSkillConfig" (full_name column), while the
# value_type column stays "Skill". Excel also emphasises the renamed
# full_name cell with an explicit (default-looking) font, which is what
# produces the extra font/style entries, and the active selection moved
# from D4 to C4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 holds the skill table's record definition:
#   B4 = full_name, C4 = value_type, D4 = define_from_file, E4 = input
# Rename the full_name from "skill.TbSkill" to "TbSkillConfig".
$ws.Range("B4").Value = "TbSkillConfig"

# Give the renamed cell its own explicit font (mirrors the new font /
# cellXfs entry that shows up in the saved styles).
$ws.Range("B4").Font.ThemeColor = 1

# The active cell/selection ends up on C4 after the edit.
[void]$ws.Range("C4").Select()
